$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.901.90'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").Value = '1.549.97'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("E4").Value = '  +0.58%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.06'
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("E6").Value = '  +1.00%  '

$ws.Range("E7").Value = '  +0.55%  '

$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.61'
$ws.Range("E9").Value = '  +1.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0588'
$ws.Range("E10").Value = '  +1.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("E11").Value = '  +0.80%  '

$ws.Range("D12").Value = '1.771.31'
$ws.Range("E12").Value = '  +0.35%  '

$ws.Range("D13").Value = '1.551.74'
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.71'
$ws.Range("E14").Value = '  +1.09%  '

$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.78'
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").Value = '26.903.77'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.90'
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.22'
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("E21").Value = '  +0.52%  '

$ws.Range("E22").Value = '  -0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.10'
$ws.Range("E23").Value = '  +1.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.24'
$ws.Range("E25").Value = '  -0.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.65'
$ws.Range("E26").Value = '  +2.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.86'
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("E28").Value = '  +0.62%  '

$ws.Range("E29").Value = '  +1.26%  '

$ws.Range("E30").Value = '  +1.10%  '

$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("E32").Value = '  +0.17%  '

$ws.Range("D33").Value = '1.405.91'
$ws.Range("E33").Value = '  +4.01%  '

$ws.Range("E34").Value = '  +2.57%  '

$ws.Range("E35").Value = '  +3.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.954'
$ws.Range("E36").Value = '  +2.23%  '

$ws.Range("E37").Value = '  +0.78%  '

$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.519'
$ws.Range("E39").Value = '  -0.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.805'
$ws.Range("E40").Value = '  +0.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.991'
$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.53'
$ws.Range("E43").Value = '  -3.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.27'
$ws.Range("E44").Value = '  +3.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.59'
$ws.Range("E45").Value = '  +1.61%  '

$ws.Range("E46").Value = '  -0.67%  '

$ws.Range("D47").Value = '1.685.36'
$ws.Range("E47").Value = '  +0.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.13'
$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("E49").Value = '  +0.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("E50").Value = '  +0.18%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0961'
$ws.Range("E51").Value = '  +3.77%  '
